$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5:AG12").ClearContents()
